$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F6").Value = -2
$ws.Range("F10").Value = -5
$ws.Range("F14").Value = -3
$ws.Range("F15").Value = -2
$ws.Range("F16").Value = 10
$ws.Range("F17").Value = 4
$ws.Range("F20").Value = -8
$ws.Range("F23").Value = 1
